$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.277776333333333
$ws.Range("H2").Value = 12.833329
$ws.Range("I2").Value = 0.3536657835996513
$ws.Range("J2").Value = 0.3568846407551645
$ws.Range("M2").Value = 0.101371
$ws.Range("N2").Value = 0.304113
$ws.Range("O2").Value = 0.0004873780190420389
$ws.Range("P2").Value = 0.0004888126654476159
$ws.Range("Q2").Value = 0.4336424646863333
$ws.Range("R2").Value = 3.902782182177
$ws.Range("S2").Value = 0.0001723689290137485
$ws.Range("T2").Value = 0.0001744497325048468

# Row 3
$ws.Range("G3").Value = 4.277776333333333
$ws.Range("H3").Value = 12.833329
$ws.Range("I3").Value = 0.3536657835996513
$ws.Range("J3").Value = 0.3568846407551645
$ws.Range("M3").Value = 123.540774
$ws.Range("N3").Value = 370.6223219999999
$ws.Range("O3").Value = 0.5939672855455723
$ws.Range("P3").Value = 0.5957156882185389
$ws.Range("Q3").Value = 528.4797992188819
$ws.Range("R3").Value = 4756.318192969937
$ws.Range("S3").Value = 0.2100659054750327
$ws.Range("T3").Value = 0.2126017793820888

# Row 4
$ws.Range("G4").Value = 4.277776333333333
$ws.Range("H4").Value = 12.833329
$ws.Range("I4").Value = 0.3536657835996513
$ws.Range("J4").Value = 0.3568846407551645
$ws.Range("M4").Value = 1.8313505
$ws.Range("N4").Value = 3.662701
$ws.Range("O4").Value = 0.008804884817764917
$ws.Range("P4").Value = 0.005887201923454927
$ws.Range("Q4").Value = 7.834107826938166
$ws.Range("R4").Value = 47.004646961629
$ws.Range("S4").Value = 0.003113986488579503
$ws.Range("T4").Value = 0.002101051943505325

# Row 5
$ws.Range("G5").Value = 4.277776333333333
$ws.Range("H5").Value = 12.833329
$ws.Range("I5").Value = 0.3536657835996513
$ws.Range("J5").Value = 0.3568846407551645
$ws.Range("M5").Value = 82.51906066666666
$ws.Range("N5").Value = 247.557182
$ws.Range("O5").Value = 0.3967404516176207
$ws.Range("P5").Value = 0.3979082971925585
$ws.Range("Q5").Value = 352.9980847687642
$ws.Range("R5").Value = 3176.982762918878
$ws.Range("S5").Value = 0.1403135227070254
$ws.Range("T5").Value = 0.1420073596970655

# Row 6
$ws.Range("G6").Value = 0.3272805
$ws.Range("H6").Value = 0.6545609999999999
$ws.Range("I6").Value = 0.02705796317293487
$ws.Range("J6").Value = 0.01820281918567982
$ws.Range("M6").Value = 0.101371
$ws.Range("N6").Value = 0.304113
$ws.Range("O6").Value = 0.0004873780190420389
$ws.Range("P6").Value = 0.0004888126654476159
$ws.Range("Q6").Value = 0.0331767515655
$ws.Range("R6").Value = 0.199060509393
$ws.Range("S6").Value = 0.00001318745649053744
$ws.Range("T6").Value = 0.000008897768564813152

# Row 7
$ws.Range("G7").Value = 0.3272805
$ws.Range("H7").Value = 0.6545609999999999
$ws.Range("I7").Value = 0.02705796317293487
$ws.Range("J7").Value = 0.01820281918567982
$ws.Range("M7").Value = 123.540774
$ws.Range("N7").Value = 370.6223219999999
$ws.Range("O7").Value = 0.5939672855455723
$ws.Range("P7").Value = 0.5957156882185389
$ws.Range("Q7").Value = 40.43248628510699
$ws.Range("R7").Value = 242.5949177106419
$ws.Range("S7").Value = 0.01607154493822019
$ws.Range("T7").Value = 0.01084370495871488

# Row 8
$ws.Range("G8").Value = 0.3272805
$ws.Range("H8").Value = 0.6545609999999999
$ws.Range("I8").Value = 0.02705796317293487
$ws.Range("J8").Value = 0.01820281918567982
$ws.Range("M8").Value = 1.8313505
$ws.Range("N8").Value = 3.662701
$ws.Range("O8").Value = 0.008804884817764917
$ws.Range("P8").Value = 0.005887201923454927
$ws.Range("Q8").Value = 0.59936530731525
$ws.Range("R8").Value = 2.397461229261
$ws.Range("S8").Value = 0.0002382422491410165
$ws.Range("T8").Value = 0.0001071636721222365

# Row 9
$ws.Range("G9").Value = 0.3272805
$ws.Range("H9").Value = 0.6545609999999999
$ws.Range("I9").Value = 0.02705796317293487
$ws.Range("J9").Value = 0.01820281918567982
$ws.Range("M9").Value = 82.51906066666666
$ws.Range("N9").Value = 247.557182
$ws.Range("O9").Value = 0.3967404516176207
$ws.Range("P9").Value = 0.3979082971925585
$ws.Range("Q9").Value = 27.006879434517
$ws.Range("R9").Value = 162.041276607102
$ws.Range("S9").Value = 0.01073498852908313
$ws.Range("T9").Value = 0.007243052786277891

# Row 10
$ws.Range("G10").Value = 7.490476666666666
$ws.Range("H10").Value = 22.47143
$ws.Range("I10").Value = 0.6192762532274139
$ws.Range("J10").Value = 0.6249125400591558
$ws.Range("M10").Value = 0.101371
$ws.Range("N10").Value = 0.304113
$ws.Range("O10").Value = 0.0004873780190420389
$ws.Range("P10").Value = 0.0004888126654476159
$ws.Range("Q10").Value = 0.7593171101766666
$ws.Range("R10").Value = 6.83385399159
$ws.Range("S10").Value = 0.000301821633537753
$ws.Range("T10").Value = 0.000305465164377956

# Row 11
$ws.Range("G11").Value = 7.490476666666666
$ws.Range("H11").Value = 22.47143
$ws.Range("I11").Value = 0.6192762532274139
$ws.Range("J11").Value = 0.6249125400591558
$ws.Range("M11").Value = 123.540774
$ws.Range("N11").Value = 370.6223219999999
$ws.Range("O11").Value = 0.5939672855455723
$ws.Range("P11").Value = 0.5957156882185389
$ws.Range("Q11").Value = 925.3792850289398
$ws.Range("R11").Value = 8328.413565260458
$ws.Range("S11").Value = 0.3678298351323195
$ws.Range("T11").Value = 0.3722702038777352

# Row 12
$ws.Range("G12").Value = 7.490476666666666
$ws.Range("H12").Value = 22.47143
$ws.Range("I12").Value = 0.6192762532274139
$ws.Range("J12").Value = 0.6249125400591558
$ws.Range("M12").Value = 1.8313505
$ws.Range("N12").Value = 3.662701
$ws.Range("O12").Value = 0.008804884817764917
$ws.Range("P12").Value = 0.005887201923454927
$ws.Range("Q12").Value = 13.71768818873833
$ws.Range("R12").Value = 82.30612913243
$ws.Range("S12").Value = 0.005452656080044398
$ws.Range("T12").Value = 0.003678986307827366

# Row 13
$ws.Range("G13").Value = 7.490476666666666
$ws.Range("H13").Value = 22.47143
$ws.Range("I13").Value = 0.6192762532274139
$ws.Range("J13").Value = 0.6249125400591558
$ws.Range("M13").Value = 82.51906066666666
$ws.Range("N13").Value = 247.557182
$ws.Range("O13").Value = 0.5939672855455723
$ws.Range("P13").Value = 0.5957156882185389
$ws.Range("Q13").Value = 528.4797992188819
$ws.Range("R13").Value = 4756.318192969937
$ws.Range("S13").Value = 0.2100659054750327
$ws.Range("T13").Value = 0.2126017793820888
